$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 2 as "bipado" (scanned) = TRUE and set the scan date/time.
$ws.Range("I2").Value = $true
$ws.Range("J2").Value = "22/05/2025 14:12"
